$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-16 Thursday", "2023-11-17 Friday"),
    @("90×74=", "74×74="),
    @("72×79=", "80×67="),
    @("17×15=", "83×18="),
    @("51×34=", "22×54="),
    @("75×27=", "15×35="),
    @("76×49=", "13×66="),
    @("52×59=", "22×74="),
    @("88×96=", "96×40="),
    @("90×87=", "44×62="),
    @("70×68=", "81×35="),
    @("74×76=", "80×93="),
    @("18×93=", "67×22="),
    @("35×27=", "81×62="),
    @("37×88=", "60×97="),
    @("53×94=", "50×52="),
    @("11×37=", "37×94="),
    @("47×93=", "69×53="),
    @("16×32=", "77×39="),
    @("21×84=", "37×36="),
    @("44×36=", "43×46="),
    @("96×14=", "39×65="),
    @("82×61=", "83×50="),
    @("68×86=", "89×41="),
    @("17×77=", "39×57="),
    @("78×36=", "78×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
